# Commit: "Fri, Jun 05, 2020  6:05:52 PM"
#
# The canonical diff shows ppt/theme/theme1.xml and ppt/theme/theme2.xml
# swapping their full contents: theme1.xml (slide master's theme, named
# "Integral") ends up holding the "Office Theme" colour scheme, and
# theme2.xml (notes master's theme, "Office Theme") ends up holding the
# "Integral" colour scheme. The two themes' font scheme and format scheme
# (fills/lines/effects) are byte-identical, so the only substantive change
# is the 10 theme colours that differ (dk2, lt2, accent1-6, hlink,
# folHlink -- dk1/lt1 are black/white in both) plus the theme/colour-scheme
# display names.
#
# This automation re-points the slide master's live theme (ppt/theme/theme1.xml)
# at the "Office Theme" palette via the PowerPoint ThemeColorScheme object
# model, which is the supported, non-file-system way to rewrite theme
# colours through COM automation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# ThemeColorScheme.Colors index order (matches msoThemeColor*�style ordering):
#   1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
#   9=accent5 10=accent6 11=hlink 12=folHlink
# RGB is the standard OLE/VBA packed colour: R + G*256 + B*65536.

$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
